$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the example data row (row 2) matching the header columns:
#   A: Document (bib key)   B: Location (page)   C: Gist/Quote   D: Q   E: Notes
$ws.Range("A2").Value = "rinker2013"
$ws.Range("B2").Value = "12"
$ws.Range("C2").Value = 'An **EXAMPLE**; feel "free" to *delete* it ***soon**'
$ws.Range("D2").Value = "y"
$ws.Range("E2").Value = "EXAMPLE (DELETE ME)"

# Move the active selection to the new data row, matching the post-edit view state.
$ws.Range("A2:E2").Select()
